$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "schwimmen"
$ws.Range("B3").Value = "dog/dog008.jpg"
$ws.Range("A4").Value = "spenden"
$ws.Range("B4").Value = "face/face004.jpg"
$ws.Range("C4").Value = "face"
$ws.Range("A6").Value = "runden"
$ws.Range("B6").Value = "face/face025.jpg"
$ws.Range("C6").Value = "face"
$ws.Range("A7").Value = "ehren"
$ws.Range("B7").Value = "dog/dog018.jpg"
$ws.Range("A9").Value = "gleichen"
$ws.Range("B9").Value = "face/face008.jpg"
$ws.Range("C9").Value = "face"
$ws.Range("A10").Value = "segnen"
$ws.Range("B10").Value = "dog/dog031.jpg"
$ws.Range("C10").Value = "dog"
$ws.Range("A12").Value = "planen"
$ws.Range("B12").Value = "dog/dog011.jpg"
$ws.Range("C12").Value = "dog"
$ws.Range("A13").Value = "schneiden"
$ws.Range("B13").Value = "dog/dog024.jpg"
$ws.Range("A15").Value = "schicken"
$ws.Range("B15").Value = "face/face018.jpg"
$ws.Range("C15").Value = "face"
$ws.Range("A16").Value = "quälen"
$ws.Range("B16").Value = "face/face024.jpg"
$ws.Range("A18").Value = "schweben"
$ws.Range("B18").Value = "face/face026.jpg"
$ws.Range("C18").Value = "face"
$ws.Range("A19").Value = "zögern"
$ws.Range("B19").Value = "dog/dog026.jpg"
$ws.Range("A21").Value = "platzen"
$ws.Range("B21").Value = "dog/dog025.jpg"
$ws.Range("A22").Value = "helfen"
$ws.Range("B22").Value = "face/face003.jpg"
$ws.Range("C22").Value = "face"
$ws.Range("A24").Value = "kleben"
$ws.Range("B24").Value = "face/face027.jpg"
$ws.Range("C24").Value = "face"
$ws.Range("A25").Value = "leugnen"
$ws.Range("B25").Value = "dog/dog017.jpg"
$ws.Range("C25").Value = "dog"
$ws.Range("A27").Value = "mühen"
$ws.Range("B27").Value = "face/face023.jpg"
$ws.Range("A28").Value = "bellen"
$ws.Range("B28").Value = "face/face015.jpg"
$ws.Range("C28").Value = "face"
$ws.Range("A30").Value = "öffnen"
$ws.Range("B30").Value = "dog/dog001.jpg"
$ws.Range("C30").Value = "dog"
$ws.Range("A31").Value = "holen"
$ws.Range("B31").Value = "dog/dog028.jpg"
$ws.Range("C31").Value = "dog"
$ws.Range("A33").Value = "faulen"
$ws.Range("B33").Value = "face/face020.jpg"
$ws.Range("A34").Value = "leisten"
$ws.Range("B34").Value = "face/face012.jpg"
$ws.Range("C34").Value = "face"
$ws.Range("A36").Value = "füttern"
$ws.Range("B36").Value = "dog/dog027.jpg"
$ws.Range("C36").Value = "dog"
$ws.Range("A37").Value = "heilen"
$ws.Range("B37").Value = "face/face005.jpg"
$ws.Range("A39").Value = "zahlen"
$ws.Range("B39").Value = "face/face007.jpg"
$ws.Range("C39").Value = "face"
$ws.Range("A40").Value = "ändern"
$ws.Range("B40").Value = "dog/dog014.jpg"
$ws.Range("C40").Value = "dog"
$ws.Range("A42").Value = "leuchten"
$ws.Range("B42").Value = "dog/dog012.jpg"
$ws.Range("C42").Value = "dog"
$ws.Range("A43").Value = "ruhen"
$ws.Range("B43").Value = "dog/dog020.jpg"
$ws.Range("C43").Value = "dog"
$ws.Range("A45").Value = "spielen"
$ws.Range("B45").Value = "face/face000.jpg"
$ws.Range("A46").Value = "testen"
$ws.Range("B46").Value = "dog/dog023.jpg"
$ws.Range("A48").Value = "werden"
$ws.Range("B48").Value = "face/face028.jpg"
$ws.Range("A49").Value = "segeln"
$ws.Range("B49").Value = "dog/dog013.jpg"
$ws.Range("C49").Value = "dog"
